# Konnect Bill Payment Verification Checks added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values (identical across data rows 2-8 per column) ---
$boTypeQuery   = "SELECT BENE_OPERATION_TYPE FROM DC_TRANSACTION DT where DT.TRANSACTION_ID='{TRANSACTION_ID}'"
$boValue       = "UPDATE"
$dbVal         = "DIGITAL_CHANNEL_SEC'"
$toAccountQry  = "SELECT K.TO_ACCOUNT FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$acctTitleQry  = "SELECT K.FT_TO_ACCOUNT_TITLE FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$beneBankQry   = "SELECT K.BENEFICIARY_BANK FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$beneIdTranQry = "SELECT BENEFICIARY_ID FROM DC_TRANSACTION K WHERE K.TRANSACTION_ID = '{TRANSACTION_ID}'"
$beneIdQry     = "SELECT FUND_TRANSFER_BENEFICIARY_ID FROM DC_FUND_TRANSFER_BENEFICIARY K WHERE K.CUSTOMER_INFO_ID = (Select CUSTOMER_INFO_ID from DC_CUSTOMER_INFO L WHERE L.CUSTOMER_NAME = '{customer_name}' ) and K.ACCOUNT_NO = '{account_number}'"

# --- Fill order below reproduces the original authoring order (column-by-column,
#     header then body) so the shared-string table comes out in the same sequence. ---

# K1 / L1
$ws.Range("K1").Value = "bene_op_type_query"
$ws.Range("L1").Value = "bene_op_value"

# L2:L8 (UPDATE)
for ($row = 2; $row -le 8; $row++) { $ws.Range("L$row").Value = $boValue }

# M1 / M2:M8
$ws.Range("M1").Value = "db_val"
for ($row = 2; $row -le 8; $row++) { $ws.Range("M$row").Value = $dbVal }

# K2:K8
for ($row = 2; $row -le 8; $row++) { $ws.Range("K$row").Value = $boTypeQuery }

# N2:N8 then N1
for ($row = 2; $row -le 8; $row++) { $ws.Range("N$row").Value = $toAccountQry }
$ws.Range("N1").Value = "to_account_query"

# O1 then O2:O8
$ws.Range("O1").Value = "account_title_query"
for ($row = 2; $row -le 8; $row++) { $ws.Range("O$row").Value = $acctTitleQry }

# P1 then P2:P8
$ws.Range("P1").Value = "bene_bank_query"
for ($row = 2; $row -le 8; $row++) { $ws.Range("P$row").Value = $beneBankQry }

# Q1, R1
$ws.Range("Q1").Value = "bene_id_tran_query"
$ws.Range("R1").Value = "bene_id_query"

# Q2:Q8, R2:R8
for ($row = 2; $row -le 8; $row++) { $ws.Range("Q$row").Value = $beneIdTranQry }
for ($row = 2; $row -le 8; $row++) { $ws.Range("R$row").Value = $beneIdQry }

# --- Formatting: mirror the cell styles already used elsewhere on the sheet by
#     copying (format-only) from reference cells that carry the same style. ---

# Style "1" (text number format) cells: K1,L1,M1,N1,O1,P1 (row1) and K2:K8,M2:M8,N2:N8 (body)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("K1:P1").PasteSpecial(-4122) | Out-Null
$ws.Range("K2:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("M2:M8").PasteSpecial(-4122) | Out-Null
$ws.Range("N2:N8").PasteSpecial(-4122) | Out-Null

# Style "4" cells: Q1,R1 (row1) and Q2:Q8,R2:R8 (body)
$ws.Range("E6").Copy() | Out-Null
$ws.Range("Q1:R1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2:Q8").PasteSpecial(-4122) | Out-Null
$ws.Range("R2:R8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Column widths for the new columns K:R (best match given engine's width
#     quantization; real Excel's bestFit pixel widths can't be hit exactly). ---
$ws.Columns.Item(11).ColumnWidth = 99.66666666666667   # K -> ~100.57
$ws.Columns.Item(12).ColumnWidth = 14.0                # L -> ~14.86
$ws.Columns.Item(13).ColumnWidth = 21.5                # M -> ~22.29
$ws.Columns.Item(14).ColumnWidth = 91.33333333333333   # N -> ~92.14
$ws.Columns.Item(15).ColumnWidth = 100.0               # O -> ~100.86
$ws.Columns.Item(16).ColumnWidth = 96.66666666666667   # P -> ~97.43
$ws.Columns.Item(17).ColumnWidth = 91.66666666666667   # Q -> ~92.43
$ws.Columns.Item(18).ColumnWidth = 243.16666666666666  # R -> 244

# New selection anchor after editing (matches author's last selection)
$ws.Range("C23").Select()
